# Lecture 9 attendance: add an "MCZ Tours" column (G) to the class roster
# and record attendance (value 2) for every student, then leave the
# selection on G13 (matches the author's final cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in G1
$ws.Range("G1").Value = "MCZ Tours"

# Attendance value of 2 for every student row (2-17)
$ws.Range("G2:G17").Value = 2

# Leave the selection where the author left it
$ws.Range("G13").Select()
